$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 73; existing rows 73-98 shift down to 74-99.
$ws.Rows.Item(73).Insert()

# Populate the newly inserted row 73 with its data (same shape as the
# surrounding "Poroto verde" records, with its own date/volume/price values).
$ws.Cells.Item(73, 1).Value = 1
$ws.Cells.Item(73, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(73, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(73, 4).Value = "2023-10-06"
$ws.Cells.Item(73, 5).Value = 15
$ws.Cells.Item(73, 6).Value = 100112031
$ws.Cells.Item(73, 7).Value = "Poroto verde"
$ws.Cells.Item(73, 8).Value = "Sin especificar"
$ws.Cells.Item(73, 9).Value = "Primera"
$ws.Cells.Item(73, 10).Value = 1800
$ws.Cells.Item(73, 11).Value = 900
$ws.Cells.Item(73, 12).Value = 1000
$ws.Cells.Item(73, 13).Value = 944
$ws.Cells.Item(73, 14).Value = "$/kilo"
$ws.Cells.Item(73, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(73, 16).Value = 944
$ws.Cells.Item(73, 17).Value = 1
$ws.Cells.Item(73, 18).Value = "Hortaliza"
